# Update Name of Algo
# Applies updated numeric results (column B and C) for specific rows
# in the result_data_KNN workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C7").Value = -13.365
$ws.Range("B9").Value = 5.468999999999999
$ws.Range("C12").Value = -11.601
$ws.Range("B13").Value = 5.601000000000001
$ws.Range("C14").Value = -12.362
$ws.Range("B16").Value = 5.294
$ws.Range("B18").Value = 5.01
$ws.Range("C19").Value = -12.028
$ws.Range("B20").Value = 7.035000000000001
$ws.Range("B26").Value = 6.223
$ws.Range("C26").Value = -13.125
$ws.Range("B27").Value = 5.680999999999999
$ws.Range("C27").Value = -13.625
$ws.Range("B29").Value = 5.373
$ws.Range("C29").Value = -11.281
$ws.Range("B35").Value = 7.657999999999999
$ws.Range("B36").Value = 7.972
$ws.Range("C37").Value = -13.313
$ws.Range("C38").Value = -13.638
$ws.Range("B45").Value = 5.601999999999999
$ws.Range("C47").Value = -12.873
$ws.Range("C51").Value = -11.282
$ws.Range("C52").Value = -11.621
$ws.Range("B55").Value = 5.354
$ws.Range("C55").Value = -13.857
$ws.Range("B57").Value = 5.363
$ws.Range("B69").Value = 5.404000000000001
$ws.Range("C69").Value = -10.732
$ws.Range("C70").Value = -11.615
$ws.Range("B76").Value = 5.98
$ws.Range("C76").Value = -12.578
$ws.Range("B78").Value = 6.904999999999999
$ws.Range("C81").Value = -13.525
$ws.Range("B82").Value = 5.061
$ws.Range("B83").Value = 5.318
$ws.Range("C83").Value = -13.668
$ws.Range("B93").Value = 5.62
$ws.Range("C94").Value = -10.869
$ws.Range("B97").Value = 5.952
$ws.Range("C100").Value = -13.409
$ws.Range("C102").Value = -13.306
